# Update resource utilization diagram with results from the latest bug fix.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Updated raw measurement data (column J = new LUTs, K = new Slice Regs,
#     L = new F7 Muxes, M = new F8 Muxes for the "Non-sharing of round keys" series) ---
$ws.Range("J3").Value = 3356
$ws.Range("J4").Value = 4603
$ws.Range("J5").Value = 5614
$ws.Range("L5").Value = 1280
$ws.Range("M5").Value = 608
$ws.Range("J6").Value = 6986
$ws.Range("L6").Value = 1536
$ws.Range("M6").Value = 736
$ws.Range("J7").Value = 8052
$ws.Range("J8").Value = 10795
$ws.Range("K8").Value = 9914
$ws.Range("J9").Value = 20514
$ws.Range("K9").Value = 17336

# Mirror the updated 15-core measurements into the summary table (row 23).
$ws.Range("C23").Value = 20514
$ws.Range("D23").Value = 17336

# --- Recompute the percentage-delta formulas (O:R) with the new throughput
#     multiplier factor ($I/$H) applied ---
$ws.Range("O3:O6").Formula = "=(J3/C3-1)*`$I3/`$H3"
$ws.Range("P3:P7").Formula = "=(K3/D3-1)*`$I3/`$H3"
$ws.Range("Q3:Q7").Formula = "=(L3/E3-1)*`$I3/`$H3"
$ws.Range("R3:R7").Formula = "=(M3/F3-1)*`$I3/`$H3"
$ws.Range("O7").Formula = "=(J7/C7-1)*`$I7/`$H7"

$ws.Range("P8:R8").Formula = "=(K8/D10-1)*`$I8/`$H8"
$ws.Range("O8").Formula = "=(J8/C10-1)*`$I8/`$H8"

$ws.Range("P9:R9").Formula = "=(K9/D17-1)*`$I9/`$H9"
$ws.Range("O9").Formula = "=(J9/C17-1)*`$I9/`$H9"

# --- Relabel the 15-core commit reference (round-keys-not-shared fix) ---
$ws.Range("B23").Value = "4a408d2: Round keys are not shared between Cores."

# --- Restore the view to the top of the sheet with M7 selected ---
$ws.Range("M7").Select()
